$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Sheet "Widgets": update household income question prompts (FR/EN) to
#    add an explanatory second line, and grow the row height to fit it.
# ---------------------------------------------------------------------------
$widgets = $wb.Worksheets.Item("Widgets")

$widgets.Range("G4").Value = "**Tranche de revenu** avant impôts (brut) **du ménage**, en 2024?`n__Cette information sert à calculer le pourcentage du revenu brut consacré aux transports et au logement.__"
$widgets.Range("H4").Value = "What was your **household's income range** before taxes (gross income), in 2024?`n__This information is used to calculate the percentage of gross income spent on transportation and housing.__"

$widgets.Rows.Item(4).RowHeight = 124.6

# ---------------------------------------------------------------------------
# 2. Sheet "Choices": bump the top two household-income brackets from
#    150-199k/200k+ to 150-209k/210k+.
# ---------------------------------------------------------------------------
$choices = $wb.Worksheets.Item("Choices")

$choices.Range("B33").Value = "150000_209999"
$choices.Range("C33").Value = "150 000$ à 209 999$"
$choices.Range("D33").Value = "$150,000 to $209,999"

$choices.Range("B34").Value = "210000_999999"
$choices.Range("C34").Value = "210 000$ et plus"
$choices.Range("D34").Value = "$210,000 and more"

# ---------------------------------------------------------------------------
# 3. Sheet "Labels": wrap the percentage-of-income label in parentheses.
# ---------------------------------------------------------------------------
$labels = $wb.Worksheets.Item("Labels")

$labels.Range("C20").Value = "({{percentageOfIncome}}% du revenu brut)"
$labels.Range("D20").Value = "({{percentageOfIncome}}% of gross income)"

# ---------------------------------------------------------------------------
# 4. Update selections on each sheet to match the author's saved cursor
#    position, then make "Choices" the active tab/sheet.
# ---------------------------------------------------------------------------
$widgets.Activate() | Out-Null
$widgets.Range("H4").Select() | Out-Null

$labels.Activate() | Out-Null
$labels.Range("C20").Select() | Out-Null

$choices.Activate() | Out-Null
$choices.Range("E32").Select() | Out-Null
